$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Caso 2")

# Title (row 2) text changed - shortened guarantee-type description
$ws2.Range("A2").Value = "Tipo Garantía Valor, Tipo Garantía Real, Tipo Garantía Fideicomiso y Tipo Garantía Aval"

# Row 5
$ws2.Range("A5").Value = "500-02-02-0900796"
$ws2.Range("C5").Value = 14450541.550000001
$ws2.Range("F5").Value = "Real-142280"
$ws2.Range("G5").Value = 0

# Row 6
$ws2.Range("A6").Value = "500-02-02-5761941"
$ws2.Range("C6").Value = 7080383.29

# Row 7 - operation label removed entirely, amount zeroed
$ws2.Range("A7").Value = $null
$ws2.Range("C7").Value = 0

# Row 9
$ws2.Range("F9").Value = "Tipo Garantía Valor, Tipo Garantía Real, Tipo Garantía Fideicomiso y Tipo Garantía Aval"

# Selection moved
$ws2.Range("C17").Select()

# Column F width adjusted (no longer auto bestFit)
$ws2.Columns("F").ColumnWidth = 12.42578125
